$wb = $excel.ActiveWorkbook

# The status text "Ready for handoff" moved to "In Translation" everywhere it
# appears (Overview!E2:F2, zh-cn!C2, de-de!C2 all share the same string).
$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus

# Those columns were sized to fit the (now shorter) status text, so their
# column widths shrink along with the text change.
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
